$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.520.64'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '3.109.09'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.58'
$ws.Range("E5").Value = '  +2.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.09'
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.440'
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.33'
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("E11").Value = '  +2.75%  '
$ws.Range("D12").Value = '3.640.95'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("E14").Value = '  +5.07%  '
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("D16").Value = '58.595.75'
$ws.Range("E16").Value = '  +1.77%  '
$ws.Range("D17").Value = '3.105.11'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.12'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.92'
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.10'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '341.24'
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.505'
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.14'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '0.0₃0921'
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("E28").Value = '  +3.86%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.27'
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.87'
$ws.Range("E31").Value = '  +3.23%  '
$ws.Range("E32").Value = '  +4.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.00'
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '154.56'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.68'
$ws.Range("E35").Value = '  +2.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.07'
$ws.Range("E36").Value = '  +3.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.31'
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("E38").Value = '  +3.98%  '
$ws.Range("D40").Value = '3.147.73'
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.680'
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.50'
$ws.Range("E44").Value = '  +8.51%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '2.288.85'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.86'
$ws.Range("E48").Value = '  +4.06%  '
$ws.Range("E49").Value = '  +2.77%  '
$ws.Range("E50").Value = '  +1.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '267.81'
$ws.Range("E51").Value = '  +6.39%  '
